$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pairs result table re-computed with additional sending/target
# clusters (ECs, FAPs, sCs) for the F2 -> F2rl1 ligand-receptor pair,
# following Dr Hou's advice. Row 1 (headers) is unchanged; rows 2-7 hold
# the full cluster x cluster combination results.

$data = @(
    @("ECs",  "F2", "F2rl1", "ECs",  3, 1,                    0.5597896666666667, 1.679369,           0.4918833646099291, 0.491883364609929,  1, 0.3333333333333333, 0.2416746666666667, 0.725024,   0.04568496189055127, 0.04568496189055127, 0.1352869810951111,  1.217582829856,     0.02247167276680075, 0.02247167276680075),
    @("ECs",  "F2", "F2rl1", "FAPs", 3, 1,                    0.5597896666666667, 1.679369,           0.4918833646099291, 0.491883364609929,  3, 1,                   5.048352,           15.145056,  0.9543150381094487,  0.9543150381094487,  2.826015283296,      25.434137549664,    0.4694116918431284,  0.4694116918431283),
    @("FAPs", "F2", "F2rl1", "ECs",  2, 0.6666666666666666,   0.4205383333333333, 1.261615,           0.3695241671379879, 0.3695241671379879, 1, 0.3333333333333333, 0.2416746666666667, 0.725024,   0.04568496189055127, 0.04568496189055127, 0.1016334615288889,  0.9147011537599999, 0.01688169749333667, 0.01688169749333667),
    @("FAPs", "F2", "F2rl1", "FAPs", 2, 0.6666666666666666,   0.4205383333333333, 1.261615,           0.3695241671379879, 0.3695241671379879, 3, 1,                   5.048352,           15.145056,  0.9543150381094487,  0.9543150381094487,  2.12302553616,       19.10722982544,     0.3526424696446512,  0.3526424696446512),
    @("sCs",  "F2", "F2rl1", "ECs",  2, 0.6666666666666666,   0.1577256666666667, 0.473177,           0.138592468252083,  0.138592468252083,  1, 0.3333333333333333, 0.2416746666666667, 0.725024,   0.04568496189055127, 0.04568496189055127, 0.03811829791644445, 0.343064681248,     0.00633159163041385, 0.006331591630413849),
    @("sCs",  "F2", "F2rl1", "FAPs", 2, 0.6666666666666666,   0.1577256666666667, 0.473177,           0.138592468252083,  0.138592468252083,  3, 1,                   5.048352,           15.145056,  0.9543150381094487,  0.9543150381094487,  0.7962546847680001,  7.166292162912001,  0.1322608766216692,  0.1322608766216691)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
